# Updated cryptos list on Tue Aug  8 06:11:39 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns with the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells hold text such as "29.222.40" or "1.000" -- mark each one as
# text (NumberFormat "@") before writing so Excel does not silently turn it
# into a number (which would drop things like trailing/leading zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.222.40"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.833.65"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "242.52"
$ws.Range("D6").Value = "0.6238"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Value = "0.07377"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").Value = "0.2927"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "23.25"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").Value = "0.07673"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "1.827.71"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").Value = "4.979"
$ws.Range("E13").Value = "  -0.72%  "
$ws.Range("D14").Value = "0.6703"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "82.71"
$ws.Range("E15").Value = "  -0.41%  "
$ws.Range("D16").Value = "0.000008954"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("D17").Value = "5.883"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "29.192.46"
$ws.Range("D19").Value = "2.074.61"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("D20").Value = "236.72"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").Value = "12.50"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").Value = "7.366"
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "158.14"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").Value = "0.1409"
$ws.Range("D27").Value = "8.552"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "17.67"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("D29").Value = "1.489"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("D30").Value = "0.05831"
$ws.Range("E30").Value = "  +4.44%  "
$ws.Range("D31").Value = "4.108"
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "4.091"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("D33").Value = "1.209"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "1.865"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "0.7334"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").Value = "2.601"
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("D38").Value = "2.855"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("D39").Value = "1.222.31"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "0.01761"
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").Value = "6.282"
$ws.Range("E41").Value = "  -3.98%  "
$ws.Range("D42").Value = "0.9067"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").Value = "101.65"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").Value = "1.981.19"
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("D46").Value = "65.56"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").Value = "0.5044"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("E48").Value = "  -3.56%  "
$ws.Range("D49").Value = "9.155"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("D50").Value = "0.4031"
$ws.Range("D51").Value = "0.1138"
$ws.Range("E51").Value = "  +3.40%  "

# Drop the text-format marker picked up above so the Price cells keep the
# workbook default style, matching their original (unstyled) appearance.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
